# Apply the "better reporting" edit: add a new captured-login row to tc_01,
# lower-case the header row on tc_02 (it now models a raw HTML form dump),
# and leave tc_03..tc_06 sample rows untouched (their data doesn't change).

$wb = $excel.ActiveWorkbook

# ---- tc_02: header text now lower-cased (new shared strings added first,
#      matching sharedStrings.xml ordering in the target workbook) ----
$ws2 = $wb.Worksheets.Item("tc_02")
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("B1").Select() | Out-Null

# ---- tc_01: new screenshot/report row, numeric login + data-provider row ----
$ws1 = $wb.Worksheets.Item("tc_01")
$ws1.Range("A2").Value = 8744954505
$ws1.Range("B2").Value = 12345
$ws1.Range("A3").Value = "q34234"
$ws1.Range("B3").Value = "dfsdf"
$ws1.Columns.Item(1).ColumnWidth = 10.14
$ws1.Range("A4:E11").Select() | Out-Null

# Restore tc_01 as the active/selected tab (selecting ranges on other sheets
# moves the active tab, so re-activate it last).
$ws1.Activate() | Out-Null

Write-Host "done"
